$d = $word.ActiveDocument

# --- Locate the two paragraphs involved in the edit -----------------------
# Paragraph 1: currently holds the "LinkedIn" hyperlink (w:hyperlink r:id="rId6")
# Paragraph 2: the (currently empty) paragraph immediately following it.
$linkedInPara = $null
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "LinkedIn*") {
        $linkedInPara = $cand
        $linkedInIndex = $i
        break
    }
}
if ($linkedInPara -eq $null) {
    throw "Could not find the 'LinkedIn' paragraph"
}
$nextPara = $d.Paragraphs.Item($linkedInIndex + 1)

# --- Replacement #1: swap the w:hyperlink run for a HYPERLINK field -------
# (begin fldChar / instrText / separate fldChar / "Portfolio Website" runs)
# The paragraph's own <w:pPr> is preserved verbatim; only its run content
# changes, matching the diff exactly.
$para1Xml = '<w:p w14:paraId="4174B704" w14:textId="0C332F41" w:rsidR="008950B0" w:rsidRDefault="009D1F5D" w:rsidP="008950B0"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:rStyle w:val="Hyperlink"/><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:fldChar w:fldCharType="begin"/></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:instrText>HYPERLINK "https://pradeep3443.github.io/"</w:instrText></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:fldChar w:fldCharType="separate"/></w:r><w:r><w:rPr><w:rStyle w:val="Hyperlink"/><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Portfoli</w:t></w:r><w:r><w:rPr><w:rStyle w:val="Hyperlink"/><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>o</w:t></w:r><w:r><w:rPr><w:rStyle w:val="Hyperlink"/><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> Website</w:t></w:r></w:p>'
$package1 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' + $para1Xml + '</w:document></pkg:xmlData></pkg:part></pkg:package>'

$linkedInPara.Range.InsertXML($package1)

# Re-resolve the second (follow-on) paragraph after the first edit, since
# indices/ranges may have shifted.
$nextPara = $d.Paragraphs.Item($linkedInIndex + 1)

# --- Replacement #2: add the matching "end" fldChar run -------------------
# Preserve the paragraph's existing <w:pPr> and simply add the new run.
$para2Xml = '<w:p w14:paraId="57768235" w14:textId="77777777" w:rsidR="008950B0" w:rsidRPr="008950B0" w:rsidRDefault="008950B0" w:rsidP="008950B0"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="467886" w:themeColor="hyperlink"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:fldChar w:fldCharType="end"/></w:r></w:p>'
$package2 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' + $para2Xml + '</w:document></pkg:xmlData></pkg:part></pkg:package>'

$nextPara.Range.InsertXML($package2)

Write-Output "Hyperlink converted to HYPERLINK field; end fldChar appended."
